$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.164.94'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.824.17'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9986'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.56'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6189'
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07346'
$ws.Range("E8").Value = '  -1.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2901'
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.02'
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07667'
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.826.38'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.957'
$ws.Range("E13").Value = '  -1.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6623'
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.16'
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008918'
$ws.Range("E16").Value = '  -4.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.834'
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.140.24'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.070.39'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.72'
$ws.Range("E20").Value = '  +6.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.41'
$ws.Range("E21").Value = '  -1.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.200'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.06'
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1417'
$ws.Range("E26").Value = '  +1.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.447'
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("E28").Value = '  -1.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.483'
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05578'
$ws.Range("E30").Value = '  -4.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.096'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.097'
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.204'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7336'
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.130'
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.622'
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.835'
$ws.Range("E38").Value = '  +2.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.218.88'
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01760'
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9211'
$ws.Range("E41").Value = '  +3.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.299'
$ws.Range("E42").Value = '  -2.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9999'
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.43'
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.973.07'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.63'
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5082'
$ws.Range("E48").Value = '  -6.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4005'
$ws.Range("E49").Value = '  -1.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.052'
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05756'
$ws.Range("E51").Value = '  -1.19%  '
